$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the rate-conversion note in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.73 = 6400.89 pesos`n✅ 6400.89 pesos = 1.71 = 925.43 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the automatic rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 579.6849999999999
$ws2.Range("O10").Value = 3710.5
$ws2.Range("N12").Value = 3735
$ws2.Range("O12").Value = 540
